$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.364.86"
$ws.Range("E2").Value = "  +4.14%  "
$ws.Range("D3").Value = "2.624.32"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'521.70"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").Value = "'141.12"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.567"
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("D9").Value = "2.623.86"
$ws.Range("E9").Value = "  +3.48%  "
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("D12").Value = "'0.331"
$ws.Range("E12").Value = "  +2.51%  "
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "3.080.50"
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").Value = "59.543.85"
$ws.Range("E15").Value = "  +4.42%  "
$ws.Range("D16").Value = "'20.44"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "2.609.00"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'338.86"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").Value = "'6.55"
$ws.Range("E22").Value = "  +7.20%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'66.47"
$ws.Range("E24").Value = "  +3.94%  "
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "0.0₃0726"
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").Value = "'5.95"
$ws.Range("E31").Value = "  -4.94%  "
$ws.Range("D32").Value = "'18.81"
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "'149.18"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "'4.01"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "'36.36"
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("D39").Value = "'0.832"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").Value = "'276.80"
$ws.Range("E42").Value = "  +7.03%  "
$ws.Range("D43").Value = "'0.997"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("E44").Value = "  +1.65%  "
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("D46").Value = "'0.0954"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0521"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'18.63"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").Value = "1.987.56"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0221"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'4.50"
$ws.Range("E51").Value = "  -0.92%  "
